$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.913.94'
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").Value = '''3.451.63'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''580.63'
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").Value = '''148.42'
$ws.Range("E6").Value = '  +0.64%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '''0.479'
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("D9").Value = '''7.91'
$ws.Range("E9").Value = '  +2.81%  '
$ws.Range("D10").Value = '''0.123'
$ws.Range("E10").Value = '  -2.47%  '
$ws.Range("D11").Value = '''0.407'
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("D12").Value = '''4.044.87'
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").Value = '''28.34'
$ws.Range("E14").Value = '  -4.60%  '
$ws.Range("D15").Value = '''3.450.87'
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("D16").Value = '''0.0000171'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").Value = '''62.979.98'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '''6.49'
$ws.Range("E18").Value = '  +2.71%  '
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").Value = '''9.15'
$ws.Range("E20").Value = '  -2.50%  '
$ws.Range("D21").Value = '''388.01'
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D22").Value = '''0.562'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '''74.81'
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '''3.594.55'
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("D26").Value = '''0.0000115'
$ws.Range("E26").Value = '  -3.18%  '
$ws.Range("D27").Value = '''0.183'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").Value = '''7.67'
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").Value = '''8.04'
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("E31").Value = '  -1.81%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").Value = '''1.35'
$ws.Range("E33").Value = '  -6.28%  '
$ws.Range("D34").Value = '''23.30'
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("D35").Value = '''1.64'
$ws.Range("E35").Value = '  +4.29%  '
$ws.Range("D36").Value = '''5.33'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '''7.04'
$ws.Range("E37").Value = '  -1.44%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '''31.66'
$ws.Range("E38").Value = '  -3.06%  '
$ws.Range("D39").Value = '''170.16'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").Value = '''3.489.24'
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("D42").Value = '''0.792'
$ws.Range("E42").Value = '  -1.69%  '
$ws.Range("D43").Value = '''42.84'
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D44").Value = '''1.71'
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("D47").Value = '''2.566.54'
$ws.Range("E47").Value = '  -2.13%  '
$ws.Range("D48").Value = '''2.27'
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("D49").Value = '''6.90'
$ws.Range("E49").Value = '  +2.14%  '
$ws.Range("D50").Value = '''22.66'
$ws.Range("E50").Value = '  -4.91%  '
$ws.Range("E51").Value = '  -0.04%  '
